# Add a new row to the "Completed" reading list for
# "The Bullet Journal Method" by Ryder Carroll.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "Completed" sheet (sheet1.xml)

# Row 67 is the current last row of data; the new entry goes in row 68.
$lastRow = 67
$newRow = $lastRow + 1

# Copy formatting (incl. the date-number-format used by columns C/D) down
# from the previous row before filling in the new values, so the new cells
# pick up the same styles (e.g. C68/D68 become short-date formatted cells
# just like C67/D67) without introducing brand-new style entries.
$ws.Range("A" + $lastRow + ":G" + $lastRow).Copy()
$ws.Range("A" + $newRow + ":G" + $newRow).PasteSpecial(-4122)

$ws.Cells.Item($newRow, 1).Value = "The Bullet Journal Method"
$ws.Cells.Item($newRow, 2).Value = "Ryder Carroll"
$ws.Cells.Item($newRow, 3).Value = 43951
$ws.Cells.Item($newRow, 4).Value = 43952
$ws.Cells.Item($newRow, 5).Value = "journaling;productivity"
$ws.Cells.Item($newRow, 6).Value = "Audio"
$ws.Cells.Item($newRow, 7).Value = "5 Hours 43 Mins"

# Reflect the post-edit UI state: the user scrolled down a bit and left the
# selection on the next empty row beneath the newly-added entry.
$null = $ws.Range("A" + ($newRow + 1)).Select()
$excel.ActiveWindow.ScrollRow = 49
